$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 78 and 79 (2025-02-20 abs_activity / rel_activity) ---
$ws.Cells.Item(78,3).Value = 7.549194711154685
$ws.Cells.Item(78,5).Value = 6.628340588899219
$ws.Cells.Item(78,6).Value = 9.642405490877502
$ws.Cells.Item(78,7).Value = 10
$ws.Cells.Item(78,10).Value = 9.100944579857059
$ws.Cells.Item(78,11).Value = 6.771123252758398
$ws.Cells.Item(78,12).Value = 8.173423997973947
$ws.Cells.Item(78,13).Value = 8.250745890327224
$ws.Cells.Item(78,16).Value = 49.19940444313953
$ws.Cells.Item(78,17).Value = 36.91677406870851

$ws.Cells.Item(79,6).Value = 6.388162183790014
$ws.Cells.Item(79,7).Value = 7.016113664601006
$ws.Cells.Item(79,10).Value = 7.166666666666666
$ws.Cells.Item(79,11).Value = 8.864247311827958
$ws.Cells.Item(79,16).Value = 33.99778299825485
$ws.Cells.Item(79,17).Value = 33.55482885045668

# --- Add new rows 82-93 (2025-02-21, 2025-02-22, 2025-02-23) ---
$ws.Cells.Item(82,1).NumberFormat = "@"
$ws.Cells.Item(82,1).Value = "2025-02-21"
$ws.Cells.Item(82,1).Style = "Normal"
$ws.Cells.Item(82,2).Value = "abs_activity"
$ws.Cells.Item(82,3).Value = 8.722749684893522
$ws.Cells.Item(82,4).Value = 0
$ws.Cells.Item(82,5).Value = 9.221475611590259
$ws.Cells.Item(82,6).Value = 9.749767230085205
$ws.Cells.Item(82,7).Value = 10
$ws.Cells.Item(82,8).Value = 10
$ws.Cells.Item(82,9).Value = 8.260408655254752
$ws.Cells.Item(82,10).Value = 5.79024057223331
$ws.Cells.Item(82,11).Value = 7.54678814283904
$ws.Cells.Item(82,12).Value = 10
$ws.Cells.Item(82,13).Value = 7.4293156017495
$ws.Cells.Item(82,14).Value = 0
$ws.Cells.Item(82,15).Value = 0
$ws.Cells.Item(82,16).Value = 51.18073769632707
$ws.Cells.Item(82,17).Value = 35.54000780231851

$ws.Cells.Item(83,1).NumberFormat = "@"
$ws.Cells.Item(83,1).Value = "2025-02-21"
$ws.Cells.Item(83,1).Style = "Normal"
$ws.Cells.Item(83,2).Value = "rel_activity"
$ws.Cells.Item(83,3).Value = 6.663956564877079
$ws.Cells.Item(83,4).Value = 5
$ws.Cells.Item(83,5).Value = 5.602772754671489
$ws.Cells.Item(83,6).Value = 6.659918801780653
$ws.Cells.Item(83,7).Value = 8.610374300029473
$ws.Cells.Item(83,8).Value = 0
$ws.Cells.Item(83,9).Value = 0
$ws.Cells.Item(83,10).Value = 0
$ws.Cells.Item(83,11).Value = 6.221875590848932
$ws.Cells.Item(83,12).Value = 0
$ws.Cells.Item(83,13).Value = 0
$ws.Cells.Item(83,14).Value = 5
$ws.Cells.Item(83,15).Value = 5
$ws.Cells.Item(83,16).Value = 32.09897921042698
$ws.Cells.Item(83,17).Value = 16.65991880178065

$ws.Cells.Item(84,1).NumberFormat = "@"
$ws.Cells.Item(84,1).Value = "2025-02-21"
$ws.Cells.Item(84,1).Style = "Normal"
$ws.Cells.Item(84,2).Value = "abs_sleep"
$ws.Cells.Item(84,3).Value = 10
$ws.Cells.Item(84,4).Value = 0
$ws.Cells.Item(84,5).Value = 10
$ws.Cells.Item(84,6).Value = 7.199999999999999
$ws.Cells.Item(84,7).Value = 10
$ws.Cells.Item(84,8).Value = 7.800000000000001
$ws.Cells.Item(84,9).Value = 8.666666666666666
$ws.Cells.Item(84,10).Value = 8.4
$ws.Cells.Item(84,11).Value = 9
$ws.Cells.Item(84,12).Value = 10
$ws.Cells.Item(84,13).Value = 0
$ws.Cells.Item(84,14).Value = 0
$ws.Cells.Item(84,15).Value = 0
$ws.Cells.Item(84,16).Value = 47.66666666666666
$ws.Cells.Item(84,17).Value = 33.4

$ws.Cells.Item(85,1).NumberFormat = "@"
$ws.Cells.Item(85,1).Value = "2025-02-21"
$ws.Cells.Item(85,1).Style = "Normal"
$ws.Cells.Item(85,2).Value = "rel_sleep"
$ws.Cells.Item(85,3).Value = 9.804893509127792
$ws.Cells.Item(85,4).Value = 0
$ws.Cells.Item(85,5).Value = 8.666699736105876
$ws.Cells.Item(85,6).Value = 0
$ws.Cells.Item(85,7).Value = 7.496557203389831
$ws.Cells.Item(85,8).Value = 8.227274368013861
$ws.Cells.Item(85,9).Value = 0
$ws.Cells.Item(85,10).Value = 0
$ws.Cells.Item(85,11).Value = 0
$ws.Cells.Item(85,12).Value = 10
$ws.Cells.Item(85,13).Value = 0
$ws.Cells.Item(85,14).Value = 0
$ws.Cells.Item(85,15).Value = 0
$ws.Cells.Item(85,16).Value = 25.9681504486235
$ws.Cells.Item(85,17).Value = 18.22727436801386

$ws.Cells.Item(86,1).NumberFormat = "@"
$ws.Cells.Item(86,1).Value = "2025-02-22"
$ws.Cells.Item(86,1).Style = "Normal"
$ws.Cells.Item(86,2).Value = "abs_activity"
$ws.Cells.Item(86,3).Value = 9.597519948473307
$ws.Cells.Item(86,4).Value = 0
$ws.Cells.Item(86,5).Value = 6.219575926922761
$ws.Cells.Item(86,6).Value = 6.672501067074402
$ws.Cells.Item(86,7).Value = 10
$ws.Cells.Item(86,8).Value = 10
$ws.Cells.Item(86,9).Value = 10
$ws.Cells.Item(86,10).Value = 10
$ws.Cells.Item(86,11).Value = 9.758160930131737
$ws.Cells.Item(86,12).Value = 2.430919834859837
$ws.Cells.Item(86,13).Value = 8.391290575406098
$ws.Cells.Item(86,14).Value = 0
$ws.Cells.Item(86,15).Value = 0
$ws.Cells.Item(86,16).Value = 53.9665473809339
$ws.Cells.Item(86,17).Value = 29.10342090193424

$ws.Cells.Item(87,1).NumberFormat = "@"
$ws.Cells.Item(87,1).Value = "2025-02-22"
$ws.Cells.Item(87,1).Style = "Normal"
$ws.Cells.Item(87,2).Value = "rel_activity"
$ws.Cells.Item(87,3).Value = 8.650139639143054
$ws.Cells.Item(87,4).Value = 5
$ws.Cells.Item(87,5).Value = 0
$ws.Cells.Item(87,6).Value = 0
$ws.Cells.Item(87,7).Value = 10
$ws.Cells.Item(87,8).Value = 0
$ws.Cells.Item(87,9).Value = 10
$ws.Cells.Item(87,10).Value = 9.67635301023622
$ws.Cells.Item(87,11).Value = 10
$ws.Cells.Item(87,12).Value = 0
$ws.Cells.Item(87,13).Value = 7.294146825396826
$ws.Cells.Item(87,14).Value = 5
$ws.Cells.Item(87,15).Value = 5
$ws.Cells.Item(87,16).Value = 50.94428646453989
$ws.Cells.Item(87,17).Value = 19.67635301023622

$ws.Cells.Item(88,1).NumberFormat = "@"
$ws.Cells.Item(88,1).Value = "2025-02-22"
$ws.Cells.Item(88,1).Style = "Normal"
$ws.Cells.Item(88,2).Value = "abs_sleep"
$ws.Cells.Item(88,3).Value = 10
$ws.Cells.Item(88,4).Value = 0
$ws.Cells.Item(88,5).Value = 10
$ws.Cells.Item(88,6).Value = 10
$ws.Cells.Item(88,7).Value = 8.966666666666667
$ws.Cells.Item(88,8).Value = 8.866666666666667
$ws.Cells.Item(88,9).Value = 10
$ws.Cells.Item(88,10).Value = 10
$ws.Cells.Item(88,11).Value = 10
$ws.Cells.Item(88,12).Value = 10
$ws.Cells.Item(88,13).Value = 10
$ws.Cells.Item(88,14).Value = 0
$ws.Cells.Item(88,15).Value = 0
$ws.Cells.Item(88,16).Value = 58.96666666666667
$ws.Cells.Item(88,17).Value = 38.86666666666667

$ws.Cells.Item(89,1).NumberFormat = "@"
$ws.Cells.Item(89,1).Value = "2025-02-22"
$ws.Cells.Item(89,1).Style = "Normal"
$ws.Cells.Item(89,2).Value = "rel_sleep"
$ws.Cells.Item(89,3).Value = 10
$ws.Cells.Item(89,4).Value = 0
$ws.Cells.Item(89,5).Value = 10
$ws.Cells.Item(89,6).Value = 10
$ws.Cells.Item(89,7).Value = 0
$ws.Cells.Item(89,8).Value = 9.274067211319799
$ws.Cells.Item(89,9).Value = 7.621047877145439
$ws.Cells.Item(89,10).Value = 10
$ws.Cells.Item(89,11).Value = 7.420877157511213
$ws.Cells.Item(89,12).Value = 10
$ws.Cells.Item(89,13).Value = 10
$ws.Cells.Item(89,14).Value = 0
$ws.Cells.Item(89,15).Value = 0
$ws.Cells.Item(89,16).Value = 45.04192503465665
$ws.Cells.Item(89,17).Value = 39.2740672113198

$ws.Cells.Item(90,1).NumberFormat = "@"
$ws.Cells.Item(90,1).Value = "2025-02-23"
$ws.Cells.Item(90,1).Style = "Normal"
$ws.Cells.Item(90,2).Value = "abs_activity"
$ws.Cells.Item(90,3).Value = 6.082292043125899
$ws.Cells.Item(90,4).Value = 0
$ws.Cells.Item(90,5).Value = 9.356597798036287
$ws.Cells.Item(90,6).Value = 0
$ws.Cells.Item(90,7).Value = 10
$ws.Cells.Item(90,8).Value = 10
$ws.Cells.Item(90,9).Value = 10
$ws.Cells.Item(90,10).Value = 10
$ws.Cells.Item(90,11).Value = -8.964677414461047
$ws.Cells.Item(90,12).Value = 1.157061367416563
$ws.Cells.Item(90,13).Value = 7.452600772450518
$ws.Cells.Item(90,14).Value = 0
$ws.Cells.Item(90,15).Value = 0
$ws.Cells.Item(90,16).Value = 33.92681319915166
$ws.Cells.Item(90,17).Value = 21.15706136741656

$ws.Cells.Item(91,1).NumberFormat = "@"
$ws.Cells.Item(91,1).Value = "2025-02-23"
$ws.Cells.Item(91,1).Style = "Normal"
$ws.Cells.Item(91,2).Value = "rel_activity"
$ws.Cells.Item(91,3).Value = 0
$ws.Cells.Item(91,4).Value = 5
$ws.Cells.Item(91,5).Value = 7.312867956265769
$ws.Cells.Item(91,6).Value = 0
$ws.Cells.Item(91,7).Value = 10
$ws.Cells.Item(91,8).Value = 10
$ws.Cells.Item(91,9).Value = 10
$ws.Cells.Item(91,10).Value = 9.511819424389536
$ws.Cells.Item(91,11).Value = 0
$ws.Cells.Item(91,12).Value = 0
$ws.Cells.Item(91,13).Value = 0
$ws.Cells.Item(91,14).Value = 5
$ws.Cells.Item(91,15).Value = 5
$ws.Cells.Item(91,16).Value = 32.31286795626577
$ws.Cells.Item(91,17).Value = 29.51181942438954

$ws.Cells.Item(92,1).NumberFormat = "@"
$ws.Cells.Item(92,1).Value = "2025-02-23"
$ws.Cells.Item(92,1).Style = "Normal"
$ws.Cells.Item(92,2).Value = "abs_sleep"
$ws.Cells.Item(92,3).Value = 10
$ws.Cells.Item(92,4).Value = 0
$ws.Cells.Item(92,5).Value = 10
$ws.Cells.Item(92,6).Value = 0
$ws.Cells.Item(92,7).Value = 10
$ws.Cells.Item(92,8).Value = 10
$ws.Cells.Item(92,9).Value = 2.133333333333333
$ws.Cells.Item(92,10).Value = 10
$ws.Cells.Item(92,11).Value = 10
$ws.Cells.Item(92,12).Value = 10
$ws.Cells.Item(92,13).Value = 10
$ws.Cells.Item(92,14).Value = 0
$ws.Cells.Item(92,15).Value = 0
$ws.Cells.Item(92,16).Value = 52.13333333333333
$ws.Cells.Item(92,17).Value = 30

$ws.Cells.Item(93,1).NumberFormat = "@"
$ws.Cells.Item(93,1).Value = "2025-02-23"
$ws.Cells.Item(93,1).Style = "Normal"
$ws.Cells.Item(93,2).Value = "rel_sleep"
$ws.Cells.Item(93,3).Value = 8.045892494929006
$ws.Cells.Item(93,4).Value = 0
$ws.Cells.Item(93,5).Value = 10
$ws.Cells.Item(93,6).Value = 0
$ws.Cells.Item(93,7).Value = 7.893802966101695
$ws.Cells.Item(93,8).Value = 10
$ws.Cells.Item(93,9).Value = 0
$ws.Cells.Item(93,10).Value = 9.995577395577399
$ws.Cells.Item(93,11).Value = 10
$ws.Cells.Item(93,12).Value = 10
$ws.Cells.Item(93,13).Value = 8.356340893677945
$ws.Cells.Item(93,14).Value = 0
$ws.Cells.Item(93,15).Value = 0
$ws.Cells.Item(93,16).Value = 44.29603635470865
$ws.Cells.Item(93,17).Value = 29.9955773955774
